# Weekly update: insert this week's two new "Coliflor" price rows
# (Primera / Segunda calidad) right before the most recent existing
# block (row 425), shifting all subsequent rows down by 2.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 2 blank rows at 425:426 - everything from the old row 425
# downward moves down by two rows (Excel's default Insert behaviour,
# which also copies formatting - e.g. the date style on column D -
# from the row above).
$ws.Rows("425:426").Insert()

# New row 425 - "Primera" calidad
$ws.Range("A425").Value = 8
$ws.Range("B425").Value = "Terminal La Palmera de La Serena"
$ws.Range("C425").Value = "Coquimbo"
$ws.Range("D425").Value = 44516
$ws.Range("E425").Value = 4
$ws.Range("F425").Value = 100112008
$ws.Range("G425").Value = "Coliflor"
$ws.Range("H425").Value = "Sin especificar"
$ws.Range("I425").Value = "Primera"
$ws.Range("J425").Value = 2300
$ws.Range("K425").Value = 550
$ws.Range("L425").Value = 600
$ws.Range("M425").Value = 575
$ws.Range("N425").Value = "$/unidad"
$ws.Range("O425").Value = "Provincia del Elquí"
$ws.Range("P425").Value = 575
$ws.Range("Q425").Value = 1
$ws.Range("R425").Value = "Hortaliza"

# New row 426 - "Segunda" calidad
$ws.Range("A426").Value = 8
$ws.Range("B426").Value = "Terminal La Palmera de La Serena"
$ws.Range("C426").Value = "Coquimbo"
$ws.Range("D426").Value = 44516
$ws.Range("E426").Value = 4
$ws.Range("F426").Value = 100112008
$ws.Range("G426").Value = "Coliflor"
$ws.Range("H426").Value = "Sin especificar"
$ws.Range("I426").Value = "Segunda"
$ws.Range("J426").Value = 1400
$ws.Range("K426").Value = 450
$ws.Range("L426").Value = 500
$ws.Range("M426").Value = 475
$ws.Range("N426").Value = "$/unidad"
$ws.Range("O426").Value = "Provincia del Elquí"
$ws.Range("P426").Value = 475
$ws.Range("Q426").Value = 1
$ws.Range("R426").Value = "Hortaliza"
